$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "63.620.11"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.74%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.600.97"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.47%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "594.58"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.00%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "150.66"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.42%  "
$ws.Range("E7").Value = "  -0.09%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.586"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("E10").Value = "  +3.36%  "
$ws.Range("E11").Value = "  +2.76%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "27.59"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.48%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "3.067.11"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.64%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "63.366.36"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.56%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000153"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +5.31%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.597.52"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.09%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "12.41"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +7.79%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.73"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +5.04%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "347.40"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("E21").Value = "  +0.35%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.19%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "67.50"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +7.93%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.32"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.89%  "
$ws.Range("E26").Value = "  -0.41%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "562.05"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.65%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "8.01"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.161"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.11%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.05"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.24%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.0₃0848"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.08%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.76"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.27%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.24"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.11%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "166.89"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.00%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.413"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +2.69%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.03%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "19.57"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +3.51%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("E40").Value = "  +0.01%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "166.93"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.78%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "39.70"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.06%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.95"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +5.49%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0587"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +4.31%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "21.99"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.29%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.630"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("E47").Value = "  +4.53%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.02"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +4.74%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0962"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.96%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "19.11"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.55%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0₆0241"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +23.00%  "
